$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells to reflect the new "Banned" terminology.
$ws.Range("B1").Value = "Banned Words"
$ws.Range("C1").Value = "Banned Websites"

# Move the active selection to D1, as left by the author after the edit.
$ws.Range("D1").Select()
